$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "0915d0e"
$ws.Range("B7").Value = "8bfdfc5"
$ws.Range("C7").Value = "1500482174.h5"
$ws.Range("D7").Value = "Model successfully drives around entire track."

$ws.Range("D7").Select()
